$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Controles")

$row = 14

$ws.Cells.Item($row, 1).Value  = "25/11/2025"       # A  Date
$ws.Cells.Item($row, 2).Value  = "19:18"            # B  Heure Prévue
$ws.Cells.Item($row, 3).Value  = "18:18"            # C  Heure Réelle
$ws.Cells.Item($row, 4).Value  = "poli"             # D  Lieu de Contrôle
$ws.Cells.Item($row, 5).Value  = "Bangoura"         # E  Nom Chauffeur
$ws.Cells.Item($row, 9).Value  = "Non observable"   # I  Type Arrêt
$ws.Cells.Item($row, 11).Value = "ras"              # K  Observation Arrêt
$ws.Cells.Item($row, 12).Value = "hombourgHaut"     # L  Client
$ws.Cells.Item($row, 24).Value = "pluvieux"         # X  Météo
$ws.Cells.Item($row, 25).Value = 22320              # Y  Parc
$ws.Cells.Item($row, 26).Value = "Conforme"         # Z  Affichage Destination
$ws.Cells.Item($row, 27).Value = "Conforme"         # AA Affichage N° Ligne
$ws.Cells.Item($row, 28).Value = "Conforme"         # AB Picto Enfant
$ws.Cells.Item($row, 29).Value = "Conforme"         # AC Tarif Affiché
$ws.Cells.Item($row, 30).Value = "Conforme"         # AD Dépliant Horaire
$ws.Cells.Item($row, 31).Value = "Conforme"         # AE Règlement
$ws.Cells.Item($row, 32).Value = "Propre"           # AF Carrosserie
$ws.Cells.Item($row, 33).Value = "ras"              # AG Observation Car
$ws.Cells.Item($row, 34).Value = "Non observable"   # AH Billetique Électronique
$ws.Cells.Item($row, 35).Value = "Conforme"         # AI Billetique Manuelle
$ws.Cells.Item($row, 36).Value = "Conforme"         # AJ Fond de Caisse
$ws.Cells.Item($row, 37).Value = "Propre"           # AK Tableau de Bord
$ws.Cells.Item($row, 38).Value = "Propre"           # AL Sol
$ws.Cells.Item($row, 39).Value = "Propre"           # AM Vitres
$ws.Cells.Item($row, 40).Value = "Propre"           # AN Sièges
$ws.Cells.Item($row, 41).Value = "ras"              # AO Observation Conditions Véhicule
$ws.Cells.Item($row, 42).Value = 10                 # AP Nombre Voyageurs
$ws.Cells.Item($row, 43).Value = 1                  # AQ Nombre Voyageurs Irréguliers
$ws.Cells.Item($row, 44).Value = "BANGOURA"         # AR Nom Contrôleur
